# new interface /& v.50.5
# Append four new ticket rows (244-247) to the tickets sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 244; A = "2024-05-23"; B = "16:04:35"; C = "No pone tornillo"; D = "-"; E = "-"; F = "-"; G = "-"; H = "16:04:37"; I = "0:00:02" },
    @{ Row = 245; A = "2024-05-23"; B = "16:04:39"; C = "Fallo tolva";      D = "-"; E = "-"; F = "-"; G = "-"; H = "16:05:06"; I = "0:00:27" },
    @{ Row = 246; A = "2024-05-23"; B = "16:24:50"; C = "-"; D = "Etiquetadora";                         E = "-"; F = "-"; G = "-"; H = "16:24:51"; I = "0:00:01" },
    @{ Row = 247; A = "2024-05-23"; B = "16:24:53"; C = "-"; D = "Detección de sealling mal puesto";     E = "-"; F = "-"; G = "-"; H = "16:24:54"; I = "0:00:01" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A holds a plain "YYYY-MM-DD" string in the source file (an
    # inline string, not a date serial). Force text storage so it isn't
    # auto-converted into a date value, then drop the temporary text
    # format so no stray cell style is left behind.
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.A
    $cellA.ClearFormats()

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
}
